$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date serial for each
# record. This was bumped by one day (45179 -> 45180) for every data row
# (rows 2 through 252) as part of an automatic refresh of the sheet.
$ws.Range("C2:C252").Value = 45180
